$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns E:F get a wider custom width, matching the existing column styling.
$ws.Range("E1:F1").ColumnWidth = 20.21875

# --- Build the header formatting once in a scratch cell, then fan it out via
# copy/paste-special so all five header cells land on a single shared style ---
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Font.Name = "宋体"
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Borders.LineStyle = 1
$scratch.Copy()

$hdr = $ws.Range("B9:F9")
$hdr.PasteSpecial(-4122)
$scratch.Clear()

$ws.Rows(9).RowHeight = 14.4

# Header labels (row 9)
$ws.Range("B9").Value = "name_recall"
$ws.Range("C9").Value = "name_precision"
$ws.Range("D9").Value = "type_accuracy"
$ws.Range("E9").Value = "value_recall"
$ws.Range("F9").Value = "value_precision"

# New recall/precision table (rows 10-15)
$ws.Range("A10").Value = "nodes_0shot_0"
$ws.Range("B10").Value = 0.73373401140391403
$ws.Range("C10").Value = 0.79257589767298475
$ws.Range("D10").Value = 0.71558021266759109
$ws.Range("E10").Value = 0.86812297734627841
$ws.Range("F10").Value = 0.89158576051779947

$ws.Range("A11").Value = "nodes_1shot_0"
$ws.Range("B11").Value = 0.77297349360456169
$ws.Range("C11").Value = 0.81281399291108025
$ws.Range("D11").Value = 0.6920711974110032
$ws.Range("E11").Value = 0.91100323624595481
$ws.Range("F11").Value = 0.92071197411003247

$ws.Range("A12").Value = "nodes_3shot_0"
$ws.Range("B12").Value = 0.83350670365233437
$ws.Range("C12").Value = 0.76533329130416505
$ws.Range("D12").Value = 0.67102018801047913
$ws.Range("E12").Value = 0.89482200647249199
$ws.Range("F12").Value = 0.89482200647249199

$ws.Range("A13").Value = "nodes_3shot_1"
$ws.Range("B13").Value = 0.84985745107104316
$ws.Range("C13").Value = 0.74554391132061015
$ws.Range("D13").Value = 0.74323470488519039
$ws.Range("E13").Value = 0.92233009708737879
$ws.Range("F13").Value = 0.91747572815533984

$ws.Range("A14").Value = "rule_cot"
$ws.Range("B14").Value = 0.80248112189859744
$ws.Range("C14").Value = 0.7348120595693407
$ws.Range("D14").Value = 0.66680536292186776
$ws.Range("E14").Value = 0.85760517799352742
$ws.Range("F14").Value = 0.88025889967637549

$ws.Range("A15").Value = "zs_cot"
$ws.Range("B15").Value = 0.59940668824163978
$ws.Range("C15").Value = 0.79729542302357803
$ws.Range("D15").Value = 0.5553513638465094
$ws.Range("E15").Value = 0.85760517799352742
$ws.Range("F15").Value = 0.87378640776699024

$ws.Range("I21").Select()
